# Added Student details in excel file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student row (row 3) beneath the existing header/data rows.
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Manoj"
$ws.Cells.Item(3, 3).Value = "Final"
$ws.Cells.Item(3, 4).Value = "JNTU"
$ws.Cells.Item(3, 5).Value = "CSE"

# Move the active selection down to A4, matching where the user left off.
$ws.Range("A4").Select()
